# "this is second commit"
# Adds a second slide (Title and Content layout) after the existing
# slide, with a title and a two-run body paragraph.

$p = $ppt.ActivePresentation

# Layout index 2 on the (single) slide master == "Title and Content"
# (ppt/slideLayouts/slideLayout2.xml) - matches the placeholders used
# by the new slide in the target deck (title + idx="1" body).
$s = $p.Slides.Add(2, 2)

# Title placeholder.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "This is second slide which is added "

# Body / content placeholder - written as two runs, matching the
# target XML ("In " + "second commit ").
$bodyRange = $s.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "In second commit "
$bodyRange.Characters(1, 3).Text = "In "
$bodyRange.Characters(4, 14).Text = "second commit "
